$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$rows = @(
    @("TestCase_B100", "OPQA-582", "Verify that more search results get displayed when user scrolls down in PATENTS search results page"),
    @("TestCase_B101", "OPQA-584", "Verify that sorting is retained when user navigates back to PATENTS search results page from record view page"),
    @("TestCase_B102", "OPQA-586", "Verify that search drop down content type is retained when user navigates back to PATENTS search results page from record view page"),
    @("TestCase_B103", "OPQA-591", "Verify that filtering is retained when user navigates back to PATENTS search results page from record view page")
)

$startRow = 101

# Copy the cell formatting (border/style) used by existing data rows so the
# new rows pick up the same cellXf (style index 3) instead of Excel minting
# brand new style/border entries.
$ws.Range("A36:E36").Copy()
$ws.Range("A101:E104").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
    $ws.Cells.Item($r, 4).Value = "Y"
}

$ws.Application.ActiveWindow.ScrollRow = 91
$ws.Range("C103").Select()
